# "Version history.xlsx" - add the 0.1.8 release row and update the
# Open points text for the 0.1.7 row (matches the upstream commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Update "Open points" (column C) for the 0.1.7 row (row 10) ---------
$ws.Range("C10").Value = "-Agrupation and disaggregate to be done in functions.`n-Change reproduction and distribution to two parts.`n-UI: Delete rows according to working functionality.`n*Implement mutations.`n*Graphic representation of F'.`n-With 4 or less niches the distribution is not equaly done.`n-Document every function.`n-Disaggregate to be done in between reproduction and distribution."

# --- Add the new 0.1.8 row (row 11) --------------------------------------
# (columns are written D before C so new shared-string entries land in the
# same order the source workbook uses: "Change log" text (47) then the
# longer "Open points" text (48))
$ws.Range("A11").Value = "0.1.8"
$ws.Range("B11").Value = "AUTOMATA CELULAR - copia (15)"
$ws.Range("D11").Value = "-Shuffle for SG corrected.`n-Output files added(_resultados.csv, _datos.csv, _nichos.csv, out.txt)."
$ws.Range("C11").Value = "-Agrupation and disaggregate to be done in functions.`n-Change reproduction and distribution to two parts.`n-UI: Delete rows according to working functionality.`n*Implement mutations.`n*Graphic representation of F'.`n-With 4 or less niches the distribution is not equaly done.`n-Document every function.`n-Disaggregate to be done in between reproduction and distribution.`n-Rework E calc."
$ws.Range("E11").Value = "Python 3.6.1"
$ws.Range("F11").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G11").Value = " PyInstaller 3.3.1"

# --- Row heights: the extra wrapped line changes the autofit height ------
$ws.Rows(10).RowHeight = 115.2
$ws.Rows(11).RowHeight = 129.6

# --- View state: keep the header frozen and focus the new row ------------
$ws.Range("A10").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C11").Select()
